$d = $word.ActiveDocument

$pairs = @(
    @("215×8=1720", "445×7=3115"),
    @("449×2=898", "453×2=906"),
    @("120×6=720", "807×3=2421"),
    @("963×8=7704", "576×7=4032"),
    @("275×7=1925", "736×2=1472"),
    @("252×4=1008", "748×7=5236"),
    @("559×8=4472", "735×6=4410"),
    @("750×4=3000", "990×5=4950"),
    @("114×5=570", "744×2=1488"),
    @("454×4=1816", "882×5=4410"),
    @("118×2=236", "438×2=876"),
    @("664×2=1328", "977×8=7816"),
    @("685×6=4110", "435×8=3480"),
    @("339×6=2034", "489×2=978"),
    @("354×4=1416", "911×9=8199"),
    @("664×9=5976", "565×6=3390"),
    @("863×7=6041", "481×4=1924"),
    @("446×6=2676", "198×8=1584"),
    @("880×5=4400", "851×9=7659"),
    @("342×2=684", "168×5=840"),
    @("949×7=6643", "784×3=2352"),
    @("659×2=1318", "707×9=6363"),
    @("831×5=4155", "954×7=6678"),
    @("485×3=1455", "677×4=2708"),
    @("638×2=1276", "558×3=1674")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
